$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 5 with the new claim record (Ramson / DL 26 K 0001)
$ws.Range("A5").Value = 9965665
$ws.Range("B5").Value = "DL 26 K 0001"
$ws.Range("C5").Value = "RAMSON"
$ws.Range("D5").Value = "65, Rajpath,, Connaught"
$ws.Range("E5").Value = "New Delhi"
$ws.Range("F5").Value = 100001
$ws.Range("G5").Value = 2147483647
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = "ramson23@gmail.com"
$ws.Range("J5").Value = "0000-00-00 00:00:00"
$ws.Range("K5").Value = "Agra"
$ws.Range("L5").Value = "Own Damage"
$ws.Range("M5").Value = "Severe Damage"
$ws.Range("N5").Value = "RAGAV"
$ws.Range("O5").Value = 60
$ws.Range("P5").Value = "TN38BXY8896668"
$ws.Range("Q5").Value = "Coimbatore North"
$ws.Range("R5").Value = "NO"
$ws.Range("S5").Value = "No"
$ws.Range("T5").Value = 43129.550312500003

# Update the active cell / selection on the sheet view
$ws.Range("A1:XFD1048576").Select()
$ws.Range("E6").Activate()
